$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Neo4j/Cypher queries stored in column B (per-tab queries) and
# column C (shared StatQuery) so they target study 'OSA01' instead of 'MGT01'.
$cells = @("B2", "B3", "B4", "B5", "C2", "C3", "C4", "C5")
foreach ($addr in $cells) {
    $range = $ws.Range($addr)
    $current = $range.Value2
    $range.Value2 = $current -replace "MGT01", "OSA01"
}

# Move the active selection from E3 to C5, matching the saved workbook view.
$ws.Range("C5").Select()
